$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 14; this pushes the
# existing rows 14-29 down to 15-30 (and keeps the sheet/table layout,
# number formats, etc. intact for the shifted rows).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 44483
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Verde"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 950
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Provincia de Linares"
$ws.Range("P14").Value = 950
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
